$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 900
$ws.Range("J10").Value = 900
$ws.Range("L10").Value = 900
$ws.Range("N10").Value = -1486

$ws.Range("H19").Value = 1517.95
$ws.Range("I19").Value = 1679.6666
$ws.Range("J19").Value = 1385.6364
$ws.Range("K19").Value = 1679.6666
$ws.Range("L19").Value = 1385.6364
$ws.Range("M19").Value = -1504.6666
$ws.Range("N19").Value = -1735.6364

$ws.Range("H40").Value = 3618.4375
$ws.Range("I40").Value = 2489.5
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 2489.5
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -2314.5
$ws.Range("N40").Value = -5850

$ws.Range("H70").Value = 17998.572
$ws.Range("I70").Value = 17998.572
$ws.Range("K70").Value = 53995.716
$ws.Range("M70").Value = -53725.716

$ws.Range("H73").Value = 17998.572
$ws.Range("I73").Value = 17998.572
$ws.Range("K73").Value = 53995.716
$ws.Range("M73").Value = -53059.716

$ws.Range("H113").Value = 6286.2856
$ws.Range("I113").Value = 2205
$ws.Range("J113").Value = 6966.5
$ws.Range("K113").Value = 2205
$ws.Range("L113").Value = 6966.5
$ws.Range("M113").Value = 1049
$ws.Range("N113").Value = -13474.5

$ws.Range("H137").Value = 11291.333
$ws.Range("I137").Value = 7001
$ws.Range("J137").Value = 12149.4
$ws.Range("K137").Value = 21003
$ws.Range("L137").Value = 36448.2
$ws.Range("M137").Value = -18453
$ws.Range("N137").Value = -41548.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1437.8
$ws.Range("I63").Value = 1547.25
$ws.Range("K63").Value = 1547.25
$ws.Range("M63").Value = -861.25

$ws.Range("H66").Value = 1437.8
$ws.Range("I66").Value = 1547.25
$ws.Range("K66").Value = 7736.25
$ws.Range("M66").Value = -4304.25

$ws.Range("H74").Value = 1914.5883
$ws.Range("I74").Value = 1581
$ws.Range("K74").Value = 1581
$ws.Range("M74").Value = -707

$ws.Range("H77").Value = 1914.5883
$ws.Range("I77").Value = 1581
$ws.Range("K77").Value = 7905
$ws.Range("M77").Value = -3537

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 22000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 22000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 22000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -22672

$ws.Range("H99").Value = 24831.275
$ws.Range("I99").Value = 1758.3478
$ws.Range("J99").Value = 113277.5
$ws.Range("K99").Value = 1758.3478
$ws.Range("L99").Value = 113277.5
$ws.Range("M99").Value = -260.3478
$ws.Range("N99").Value = -116273.5

$ws.Range("H132").Value = 152523
$ws.Range("J132").Value = 152523
$ws.Range("L132").Value = 152523
$ws.Range("N132").Value = -162643

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 34417.5
$ws.Range("I5").Value = 168.66667
$ws.Range("J5").Value = 68666.336
$ws.Range("K5").Value = 168.66667
$ws.Range("L5").Value = 68666.336
$ws.Range("M5").Value = -56.66667000000001
$ws.Range("N5").Value = -68890.336

$ws.Range("H31").Value = 3269.7856
$ws.Range("I31").Value = 2198.5312
$ws.Range("J31").Value = 6697.8
$ws.Range("K31").Value = 2198.5312
$ws.Range("L31").Value = 6697.8
$ws.Range("M31").Value = -1903.5312
$ws.Range("N31").Value = -7287.8

$ws.Range("H34").Value = 3269.7856
$ws.Range("I34").Value = 2198.5312
$ws.Range("J34").Value = 6697.8
$ws.Range("K34").Value = 2198.5312
$ws.Range("L34").Value = 6697.8
$ws.Range("M34").Value = -1996.5312
$ws.Range("N34").Value = -7101.8

$ws.Range("H39").Value = 850.4
$ws.Range("I39").Value = 850.4
$ws.Range("K39").Value = 850.4
$ws.Range("M39").Value = -459.4

$ws.Range("H43").Value = 26964.5
$ws.Range("J43").Value = 26964.5
$ws.Range("L43").Value = 26964.5
$ws.Range("N43").Value = -27332.5

$ws.Range("H49").Value = 850.4
$ws.Range("I49").Value = 850.4
$ws.Range("K49").Value = 850.4
$ws.Range("M49").Value = -668.4

$ws.Range("H86").Value = 71433900
$ws.Range("I86").Value = 166670580
$ws.Range("J86").Value = 6397.375
$ws.Range("K86").Value = 166670580
$ws.Range("L86").Value = 6397.375
$ws.Range("M86").Value = -166669457
$ws.Range("N86").Value = -8643.375

$ws.Range("H89").Value = 71433900
$ws.Range("I89").Value = 166670580
$ws.Range("J89").Value = 6397.375
$ws.Range("K89").Value = 833352900
$ws.Range("L89").Value = 31986.875
$ws.Range("M89").Value = -833347284
$ws.Range("N89").Value = -43218.875

$ws.Range("H96").Value = 18493.125
$ws.Range("J96").Value = 18493.125
$ws.Range("L96").Value = 18493.125
$ws.Range("N96").Value = -23985.125

$ws.Range("H99").Value = 31641.438
$ws.Range("I99").Value = 60271.5
$ws.Range("J99").Value = 3011.375
$ws.Range("K99").Value = 60271.5
$ws.Range("L99").Value = 3011.375
$ws.Range("M99").Value = -58773.5
$ws.Range("N99").Value = -6007.375

$ws.Range("H101").Value = 26964.5
$ws.Range("J101").Value = 26964.5
$ws.Range("L101").Value = 26964.5
$ws.Range("N101").Value = -33454.5

$ws.Range("H126").Value = 31641.438
$ws.Range("I126").Value = 60271.5
$ws.Range("J126").Value = 3011.375
$ws.Range("K126").Value = 180814.5
$ws.Range("L126").Value = 9034.125
$ws.Range("M126").Value = -178344.5
$ws.Range("N126").Value = -13974.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2288.25
$ws.Range("J35").Value = 4001.5
$ws.Range("L35").Value = 12004.5
$ws.Range("N35").Value = -12580.5

$ws.Range("H58").Value = 7163.375
$ws.Range("J58").Value = 7635.3335
$ws.Range("L58").Value = 22906.0005
$ws.Range("N58").Value = -23162.0005

$ws.Range("H122").Value = 271.2857
$ws.Range("I122").Value = 137.5
$ws.Range("J122").Value = 324.8
$ws.Range("K122").Value = 1237.5
$ws.Range("L122").Value = 2923.2
$ws.Range("M122").Value = 1212.5
$ws.Range("N122").Value = -7823.200000000001

$ws.Range("H132").Value = 2074.1875
$ws.Range("J132").Value = 2074.1875
$ws.Range("L132").Value = 18667.6875
$ws.Range("N132").Value = -23727.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 32554.25
$ws.Range("I43").Value = 200
$ws.Range("J43").Value = 43339
$ws.Range("K43").Value = 200
$ws.Range("L43").Value = 43339
$ws.Range("M43").Value = -49
$ws.Range("N43").Value = -43641

$ws.Range("H80").Value = 111124550
$ws.Range("I80").Value = 111124550
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 111124550
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -111123552
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 111124550
$ws.Range("I83").Value = 111124550
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 555622750
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -555617758
$ws.Range("N83").ClearContents()

$ws.Range("H126").Value = 48432
$ws.Range("I126").Value = 1728
$ws.Range("J126").Value = 110704
$ws.Range("K126").Value = 5184
$ws.Range("L126").Value = 332112
$ws.Range("M126").Value = -2714
$ws.Range("N126").Value = -337052

$ws.Range("H132").Value = 5215.5557
$ws.Range("J132").Value = 4422.4546
$ws.Range("L132").Value = 13267.3638
$ws.Range("N132").Value = -18327.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25285.285
$ws.Range("I7").Value = 28499.5
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 28499.5
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -28387.5
$ws.Range("N7").Value = -6224

$ws.Range("H61").Value = 15891908
$ws.Range("I61").Value = 19610346
$ws.Range("K61").Value = 19610346
$ws.Range("M61").Value = -19610144

$ws.Range("H103").Value = 69000
$ws.Range("J103").Value = 69000
$ws.Range("L103").Value = 69000
$ws.Range("N103").Value = -71344

$ws.Range("H113").Value = 15891908
$ws.Range("I113").Value = 19610346
$ws.Range("K113").Value = 19610346
$ws.Range("M113").Value = -19608176

$ws.Range("H126").Value = 25285.285
$ws.Range("I126").Value = 28499.5
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 85498.5
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -83028.5
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 4126.2
$ws.Range("I132").Value = 3579.5386
$ws.Range("K132").Value = 10738.6158
$ws.Range("M132").Value = -8208.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 37557
$ws.Range("J64").Value = 37557
$ws.Range("L64").Value = 37557
$ws.Range("N64").Value = -38053

$ws.Range("H67").Value = 37557
$ws.Range("J67").Value = 37557
$ws.Range("L67").Value = 37557
$ws.Range("N67").Value = -39273

$ws.Range("H81").Value = 250001600
$ws.Range("I81").Value = 333334300
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 666668600
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -666667539
$ws.Range("N81").Value = -9122

$ws.Range("H84").Value = 250001600
$ws.Range("I84").Value = 333334300
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 3333343000
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -3333337696
$ws.Range("N84").Value = -45608
